$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.340.73'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.486.08'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '568.54'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.64'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('D9').Value = '2.483.14'
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.158'
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.86'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').Value = '2.940.19'
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('D15').Value = '69.174.17'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '24.03'
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('D18').Value = '2.480.56'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.13'
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.35'
$ws.Range('E20').Value = '  -4.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '345.38'
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.88'
$ws.Range('E24').Value = '  -5.58%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '69.17'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.87'
$ws.Range('E26').Value = '  -3.30%  '
$ws.Range('D27').Value = '2.621.44'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.57'
$ws.Range('E28').Value = '  -4.11%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.02'
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('D30').Value = '0.0₃0861'
$ws.Range('E30').Value = '  -3.86%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.53'
$ws.Range('E31').Value = '  -4.68%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '435.34'
$ws.Range('E32').Value = '  -6.36%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '156.74'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.113'
$ws.Range('E37').Value = '  -3.33%  '
$ws.Range('B38').Value = 'WhiteBITCoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.05'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.06'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.312'
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('E42').Value = '  -4.74%  '
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('B44').Value = 'POPCAT'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.28'
$ws.Range('E44').Value = '  +42.07%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.14'
$ws.Range('E45').Value = '  -6.30%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.06'
$ws.Range('E46').Value = '  -6.50%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '137.69'
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.40'
$ws.Range('E48').Value = '  -2.41%  '
$ws.Range('E49').Value = '  -4.62%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0720'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.571'
$ws.Range('E51').Value = '  -1.14%  '
